$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$rangeA = $ws.Range("A1:A$lastRow")
$rangeB = $ws.Range("B1:B$lastRow")

$valuesA = $rangeA.Value2
$valuesB = $rangeB.Value2

$rangeA.Value2 = $valuesB
$rangeB.Value2 = $valuesA
